$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.120.19"
$ws.Range("E2").Value = "  +5.57%  "
$ws.Range("D3").Value = "2.495.98"
$ws.Range("E3").Value = "  +4.12%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'161.92"
$ws.Range("E5").Value = "  +9.54%  "
$ws.Range("D6").Value = "'503.80"
$ws.Range("E6").Value = "  +5.16%  "
$ws.Range("D7").Value = "'0.613"
$ws.Range("E7").Value = "  +22.77%  "
$ws.Range("D8").Value = "'0.991"
$ws.Range("E8").Value = "  -0.81%  "
$ws.Range("D9").Value = "2.528.70"
$ws.Range("E9").Value = "  +5.06%  "
$ws.Range("E10").Value = "  +15.62%  "
$ws.Range("D11").Value = "'0.104"
$ws.Range("E11").Value = "  +7.10%  "
$ws.Range("D12").Value = "'0.341"
$ws.Range("E12").Value = "  +5.73%  "
$ws.Range("D13").Value = "'0.127"
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("D14").Value = "2.906.64"
$ws.Range("E14").Value = "  +3.24%  "
$ws.Range("D15").Value = "58.980.18"
$ws.Range("E15").Value = "  +4.73%  "
$ws.Range("D16").Value = "'22.08"
$ws.Range("E16").Value = "  +8.53%  "
$ws.Range("D17").Value = "'0.0000138"
$ws.Range("E17").Value = "  +4.59%  "
$ws.Range("D18").Value = "2.508.16"
$ws.Range("E18").Value = "  +4.37%  "
$ws.Range("D19").Value = "'4.75"
$ws.Range("E19").Value = "  +5.76%  "
$ws.Range("D20").Value = "'333.69"
$ws.Range("E20").Value = "  +5.86%  "
$ws.Range("D21").Value = "'10.24"
$ws.Range("E21").Value = "  +4.68%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'6.02"
$ws.Range("E22").Value = "  +5.92%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").Value = "'59.29"
$ws.Range("E24").Value = "  +4.11%  "
$ws.Range("D25").Value = "'0.415"
$ws.Range("E25").Value = "  +4.86%  "
$ws.Range("D26").Value = "'0.167"
$ws.Range("E26").Value = "  +5.94%  "
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").Value = "2.572.18"
$ws.Range("E28").Value = "  +2.75%  "
$ws.Range("D29").Value = "'7.53"
$ws.Range("E29").Value = "  +2.98%  "
$ws.Range("D30").Value = "0.0₃0818"
$ws.Range("E30").Value = "  +5.77%  "
$ws.Range("D31").Value = "'0.996"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").Value = "'19.30"
$ws.Range("E32").Value = "  +7.32%  "
$ws.Range("D33").Value = "'154.47"
$ws.Range("E33").Value = "  +4.02%  "
$ws.Range("D34").Value = "'1.56"
$ws.Range("E34").Value = "  +5.03%  "
$ws.Range("D35").Value = "'5.52"
$ws.Range("E35").Value = "  +8.18%  "
$ws.Range("D36").Value = "'3.93"
$ws.Range("E36").Value = "  +9.37%  "
$ws.Range("D37").Value = "'1.19"
$ws.Range("E37").Value = "  +7.46%  "
$ws.Range("D38").Value = "'0.862"
$ws.Range("E38").Value = "  +2.42%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "'3.72"
$ws.Range("E39").Value = "  +10.00%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'1.44"
$ws.Range("E40").Value = "  +7.19%  "
$ws.Range("D41").Value = "'34.87"
$ws.Range("E41").Value = "  +4.12%  "
$ws.Range("D42").Value = "'288.92"
$ws.Range("E42").Value = "  +13.74%  "
$ws.Range("E43").Value = "  +7.20%  "
$ws.Range("D44").Value = "'0.617"
$ws.Range("E44").Value = "  +5.27%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0555"
$ws.Range("E45").Value = "  +2.34%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "'0.991"
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("D47").Value = "'0.0238"
$ws.Range("E47").Value = "  +6.63%  "
$ws.Range("D48").Value = "'4.81"
$ws.Range("E48").Value = "  +3.63%  "
$ws.Range("D49").Value = "'10.27"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("D50").Value = "'0.716"
$ws.Range("E50").Value = "  +13.55%  "
$ws.Range("D51").Value = "'18.39"
$ws.Range("E51").Value = "  +8.02%  "

# Reset style on cells where a quote-prefix was used to force text,
# so no stray number-format / quotePrefix style remains on the cell.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
